# Add a default header to the section containing the questionnaire number,
# e.g. "Questionnaire 34", so printed copies can be tracked.

$d = $word.ActiveDocument

$section = $d.Sections(1)
$header = $section.Headers(1)          # wdHeaderFooterPrimary

$r = $header.Range
$r.InsertAfter("Questionnaire 34")
$r.Style = "Header"
$r.ParagraphFormat.Alignment = 1       # wdAlignParagraphCenter

# Apply the run-level font formatting only to the inserted text, not to the
# trailing paragraph mark, by excluding the last character of the range.
$runRange = $r.Duplicate
[void]$runRange.MoveEnd(1, -1)
$runRange.Font.Name = "Arial"
$runRange.Font.Size = 12
